# Yang Li is trying to add his name again as the previous attempt seems not working.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of data (Name / Completed).
$ws.Range("A3").Value = "Yang Li"
$ws.Range("B3").Value = "Yes"

# Match the page setup that was applied when the sheet was printed/saved.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where Excel would land after typing the row (next row, col A).
$ws.Range("A4").Select() | Out-Null
